# Auto-generated Excel COM-interop script
# Applies cell value updates per the commit diff (scheduled-runner price refresh)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 61.8
$ws.Range("I8").Value = 61.8
$ws.Range("K8").Value = 185.4
$ws.Range("M8").Value = -46.39999999999998
$ws.Range("H32").Value = 884.5714
$ws.Range("I32").Value = 600
$ws.Range("J32").Value = 998.4
$ws.Range("K32").Value = 600
$ws.Range("L32").Value = 998.4
$ws.Range("M32").Value = -274
$ws.Range("N32").Value = -1650.4
$ws.Range("H64").Value = 3821.4285
$ws.Range("I64").Value = 3675
$ws.Range("J64").Value = 4016.6667
$ws.Range("K64").Value = 3675
$ws.Range("L64").Value = 4016.6667
$ws.Range("M64").Value = -3427
$ws.Range("N64").Value = -4512.6667
$ws.Range("H67").Value = 3821.4285
$ws.Range("I67").Value = 3675
$ws.Range("J67").Value = 4016.6667
$ws.Range("K67").Value = 3675
$ws.Range("L67").Value = 4016.6667
$ws.Range("M67").Value = -2817
$ws.Range("N67").Value = -5732.6667
$ws.Range("H74").Value = 4166.222
$ws.Range("I74").Value = 3999
$ws.Range("K74").Value = 3999
$ws.Range("M74").Value = -3063
$ws.Range("H76").Value = 3048
$ws.Range("I76").Value = 3036.3635
$ws.Range("J76").Value = 3133.3333
$ws.Range("K76").Value = 3036.3635
$ws.Range("L76").Value = 3133.3333
$ws.Range("M76").Value = -2721.3635
$ws.Range("N76").Value = -3763.3333
$ws.Range("H77").Value = 4166.222
$ws.Range("I77").Value = 3999
$ws.Range("K77").Value = 19995
$ws.Range("M77").Value = -15315
$ws.Range("H79").Value = 3048
$ws.Range("I79").Value = 3036.3635
$ws.Range("J79").Value = 3133.3333
$ws.Range("K79").Value = 3036.3635
$ws.Range("L79").Value = 3133.3333
$ws.Range("M79").Value = -1944.3635
$ws.Range("N79").Value = -5317.3333
$ws.Range("H132").Value = 2443.276
$ws.Range("I132").Value = 1805.9615
$ws.Range("J132").Value = 7966.6665
$ws.Range("K132").Value = 5417.8845
$ws.Range("L132").Value = 23899.9995
$ws.Range("M132").Value = -2887.8845
$ws.Range("N132").Value = -28959.9995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3099.5574
$ws.Range("I32").Value = 2574.0544
$ws.Range("J32").Value = 7916.6665
$ws.Range("K32").Value = 2574.0544
$ws.Range("L32").Value = 7916.6665
$ws.Range("M32").Value = -2287.0544
$ws.Range("N32").Value = -8490.666499999999
$ws.Range("H45").Value = 1491.7391
$ws.Range("I45").Value = 1051.8334
$ws.Range("J45").Value = 1971.6364
$ws.Range("K45").Value = 1051.8334
$ws.Range("L45").Value = 1971.6364
$ws.Range("M45").Value = -674.8334
$ws.Range("N45").Value = -2725.6364
$ws.Range("H88").Value = 2422
$ws.Range("I88").Value = 2189.75
$ws.Range("J88").Value = 2590.9092
$ws.Range("K88").Value = 2189.75
$ws.Range("L88").Value = 2590.9092
$ws.Range("M88").Value = -1783.75
$ws.Range("N88").Value = -3402.9092
$ws.Range("H91").Value = 2422
$ws.Range("I91").Value = 2189.75
$ws.Range("J91").Value = 2590.9092
$ws.Range("K91").Value = 2189.75
$ws.Range("L91").Value = 2590.9092
$ws.Range("M91").Value = -785.75
$ws.Range("N91").Value = -5398.9092

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H57").Value = 48000
$ws.Range("J57").Value = 48000
$ws.Range("L57").Value = 48000
$ws.Range("N57").Value = -49440
$ws.Range("H86").Value = 1747.4
$ws.Range("I86").Value = 1631.6154
$ws.Range("J86").Value = 2500
$ws.Range("K86").Value = 1631.6154
$ws.Range("L86").Value = 2500
$ws.Range("M86").Value = -508.6153999999999
$ws.Range("N86").Value = -4746
$ws.Range("H89").Value = 1747.4
$ws.Range("I89").Value = 1631.6154
$ws.Range("J89").Value = 2500
$ws.Range("K89").Value = 8158.076999999999
$ws.Range("L89").Value = 12500
$ws.Range("M89").Value = -2542.076999999999
$ws.Range("N89").Value = -23732
$ws.Range("H105").Value = 2212.2173
$ws.Range("I105").Value = 1841.1111
$ws.Range("J105").Value = 2450.7856
$ws.Range("K105").Value = 1841.1111
$ws.Range("L105").Value = 2450.7856
$ws.Range("M105").Value = -94.11110000000008
$ws.Range("N105").Value = -5944.7856
$ws.Range("H136").Value = 48000
$ws.Range("J136").Value = 48000
$ws.Range("L136").Value = 48000
$ws.Range("N136").Value = -58200

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 1994.75
$ws.Range("I6").Value = 1000
$ws.Range("K6").Value = 1000
$ws.Range("M6").Value = -887
$ws.Range("H25").Value = 8980
$ws.Range("I25").Value = 7450
$ws.Range("J25").Value = 10000
$ws.Range("K25").Value = 7450
$ws.Range("L25").Value = 10000
$ws.Range("M25").Value = -7276
$ws.Range("H62").Value = 3240.6
$ws.Range("J62").Value = 4626.5
$ws.Range("L62").Value = 4626.5
$ws.Range("N62").Value = -5874.5
$ws.Range("H65").Value = 3240.6
$ws.Range("J65").Value = 4626.5
$ws.Range("L65").Value = 23132.5
$ws.Range("N65").Value = -29372.5
$ws.Range("N25").Value = -10348

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 3516.25
$ws.Range("I17").Value = 157
$ws.Range("J17").Value = 6875.5
$ws.Range("K17").Value = 471
$ws.Range("L17").Value = 20626.5
$ws.Range("M17").Value = -302
$ws.Range("N17").Value = -20964.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2904.2222
$ws.Range("J80").Value = 3512.6667
$ws.Range("L80").Value = 3512.6667
$ws.Range("N80").Value = -5508.6667
$ws.Range("H83").Value = 2904.2222
$ws.Range("J83").Value = 3512.6667
$ws.Range("L83").Value = 17563.3335
$ws.Range("N83").Value = -27547.3335

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3150.25
$ws.Range("I7").Value = 2867.3333
$ws.Range("J7").Value = 3999
$ws.Range("K7").Value = 2867.3333
$ws.Range("L7").Value = 3999
$ws.Range("M7").Value = -2755.3333
$ws.Range("N7").Value = -4223
$ws.Range("H22").Value = 100613
$ws.Range("I22").Value = 333766.66
$ws.Range("J22").Value = 690
$ws.Range("K22").Value = 333766.66
$ws.Range("L22").Value = 690
$ws.Range("N22").Value = -1280
$ws.Range("H27").Value = 100613
$ws.Range("I27").Value = 333766.66
$ws.Range("J27").Value = 690
$ws.Range("K27").Value = 333766.66
$ws.Range("L27").Value = 690
$ws.Range("N27").Value = -904
$ws.Range("H61").Value = 3625
$ws.Range("I61").Value = 1500
$ws.Range("K61").Value = 1500
$ws.Range("M61").Value = -1298
$ws.Range("H68").Value = 6865.087
$ws.Range("I68").Value = 8766.267
$ws.Range("J68").Value = 3300.375
$ws.Range("K68").Value = 8766.267
$ws.Range("L68").Value = 3300.375
$ws.Range("M68").Value = -8017.267
$ws.Range("N68").Value = -4798.375
$ws.Range("H71").Value = 6865.087
$ws.Range("I71").Value = 8766.267
$ws.Range("J71").Value = 3300.375
$ws.Range("K71").Value = 43831.335
$ws.Range("L71").Value = 16501.875
$ws.Range("M71").Value = -40087.335
$ws.Range("N71").Value = -23989.875
$ws.Range("H82").Value = 2201.25
$ws.Range("I82").Value = 650
$ws.Range("J82").Value = 2422.8572
$ws.Range("K82").Value = 650
$ws.Range("L82").Value = 2422.8572
$ws.Range("M82").Value = -289
$ws.Range("N82").Value = -3144.8572
$ws.Range("H85").Value = 2201.25
$ws.Range("I85").Value = 650
$ws.Range("J85").Value = 2422.8572
$ws.Range("K85").Value = 650
$ws.Range("L85").Value = 2422.8572
$ws.Range("M85").Value = 598
$ws.Range("N85").Value = -4918.8572
$ws.Range("H113").Value = 3625
$ws.Range("I113").Value = 1500
$ws.Range("K113").Value = 1500
$ws.Range("M113").Value = 670
$ws.Range("H126").Value = 3150.25
$ws.Range("I126").Value = 2867.3333
$ws.Range("J126").Value = 3999
$ws.Range("K126").Value = 8601.999899999999
$ws.Range("L126").Value = 11997
$ws.Range("M126").Value = -6131.999899999999
$ws.Range("N126").Value = -16937
$ws.Range("M22").Value = -333471.66
$ws.Range("M27").Value = -333659.66
